# Applies "added harvard case classification" update:
#  - swaps the average_doctor / average_doctor_old header labels (BP1/BQ1)
#  - recomputes the precision/recall/f1/f2/NDCG stats (rows 4-8) for the
#    "_old" app columns and the average_doctor / average_doctor_old columns
#  - shifts the average_doctor(_old) values for the remaining stat rows (9-13)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("AI4").Value = 0.1
$ws.Range("AJ4").Value = 0.023
$ws.Range("AK4").Value = 0.153
$ws.Range("BA4").Value = 1.624
$ws.Range("BB4").Value = 0.175
$ws.Range("BC4").Value = 0.419
$ws.Range("BG4").Value = 0.583
$ws.Range("BH4").Value = 0.201
$ws.Range("BI4").Value = 0.449
$ws.Range("BP4").Value = 0.541
$ws.Range("BQ4").Value = 0.576
$ws.Range("E4").Value = 0.617
$ws.Range("F4").Value = 0.08500000000000001
$ws.Range("G4").Value = 0.291
$ws.Range("N4").Value = 0.472
$ws.Range("O4").Value = 0.06900000000000001
$ws.Range("P4").Value = 0.262
$ws.Range("W4").Value = 0.167
$ws.Range("X4").Value = 0.065
$ws.Range("Y4").Value = 0.255
$ws.Range("AI5").Value = 0.139
$ws.Range("AJ5").Value = 0.041
$ws.Range("AK5").Value = 0.202
$ws.Range("BA5").Value = 1.042
$ws.Range("BB5").Value = 0.05
$ws.Range("BC5").Value = 0.224
$ws.Range("BG5").Value = 0.361
$ws.Range("BH5").Value = 0.06900000000000001
$ws.Range("BI5").Value = 0.262
$ws.Range("BP5").Value = 0.347
$ws.Range("BQ5").Value = 0.373
$ws.Range("E5").Value = 0.514
$ws.Range("F5").Value = 0.024
$ws.Range("G5").Value = 0.155
$ws.Range("N5").Value = 0.708
$ws.Range("O5").Value = 0.06
$ws.Range("P5").Value = 0.244
$ws.Range("W5").Value = 0.139
$ws.Range("X5").Value = 0.041
$ws.Range("Y5").Value = 0.202
$ws.Range("AI6").Value = 0.116
$ws.Range("BA6").Value = 1.269
$ws.Range("BG6").Value = 0.446
$ws.Range("BP6").Value = 0.423
$ws.Range("BQ6").Value = 0.452
$ws.Range("E6").Value = 0.5610000000000001
$ws.Range("N6").Value = 0.5659999999999999
$ws.Range("W6").Value = 0.152
$ws.Range("AI7").Value = 0.129
$ws.Range("BA7").Value = 1.123
$ws.Range("BG7").Value = 0.391
$ws.Range("BP7").Value = 0.374
$ws.Range("BQ7").Value = 0.401
$ws.Range("E7").Value = 0.532
$ws.Range("N7").Value = 0.644
$ws.Range("W7").Value = 0.144
$ws.Range("AI8").Value = 0.06
$ws.Range("AJ8").Value = 0.008
$ws.Range("AK8").Value = 0.089
$ws.Range("BA8").Value = 1.437
$ws.Range("BB8").Value = 0.125
$ws.Range("BC8").Value = 0.354
$ws.Range("BG8").Value = 0.473
$ws.Range("BH8").Value = 0.156
$ws.Range("BI8").Value = 0.395
$ws.Range("BP8").Value = 0.479
$ws.Range("BQ8").Value = 0.511
$ws.Range("E8").Value = 0.546
$ws.Range("F8").Value = 0.061
$ws.Range("G8").Value = 0.247
$ws.Range("N8").Value = 0.902
$ws.Range("O8").Value = 0.008999999999999999
$ws.Range("P8").Value = 0.093
$ws.Range("W8").Value = 0.184
$ws.Range("X8").Value = 0.08
$ws.Range("Y8").Value = 0.283
$ws.Range("BP10").Value = 0.555
$ws.Range("BQ10").Value = 0.648
$ws.Range("BP11").Value = 0.555
$ws.Range("BQ11").Value = 0.648
$ws.Range("BP12").Value = 1.067
$ws.Range("BQ12").Value = 1.353
$ws.Range("BP13").Value = 0.732
$ws.Range("BQ13").Value = 0.739
